$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "HTSE": insert a new "kg-H2/kWh-e" row (inverse of electricity
# required) right under "Electricity required", highlight the NOAK
# numbers in green ("Good" style) and add a footnote explaining the
# NOAK values were used.
# ------------------------------------------------------------------
$htse = $wb.Worksheets.Item("HTSE")

# Insert a new row 4 (pushes everything from the old row 4 down by one)
$htse.Rows.Item(4).Insert()

$htse.Range("B4").Formula = "=1/B3"
$htse.Range("C4").Value = "kg-H2/kWh-e"

# Highlight the NOAK-sourced numbers with the built-in "Good" (green) style
$htse.Range("B3").Style = "Good"
$htse.Range("B4").Style = "Good"
$htse.Range("B9").Style = "Good"
$htse.Range("B11").Style = "Good"
$htse.Range("B12").Style = "Good"
$htse.Range("B14").Style = "Good"

# Footnote row
$htse.Range("A16").Value = "Using NOAK values"
$htse.Range("A16").Style = "Good"

# Column width tweaks
$htse.Columns.Item(1).ColumnWidth = 29
$htse.Columns.Item(2).ColumnWidth = 15.36328125
$htse.Columns.Item(3).ColumnWidth = 12.08984375

$htse.Range("H14").Select()

# ------------------------------------------------------------------
# Sheet "Transfer_rates": the electricity line of the FT process is
# fixed no matter the production level (it is modelled as its own
# component), so drop the stale numbers and annotate instead. Also
# bold/merge the two section header rows (row 1 "FT process" and
# row 11 "HTSE") like the existing "Inputs"/"Outputs" header rows,
# and tint the unit-conversion helper columns gray.
# ------------------------------------------------------------------
$tr = $wb.Worksheets.Item("Transfer_rates")

$tr.Range("A1:G1").Style = "Normal"
$tr.Range("A1:G1").Font.Bold = $true
$tr.Range("A1:G1").Merge()

$tr.Range("A11:G11").Font.Bold = $true
$tr.Range("A11:G11").Merge()

$tr.Range("D3").Font.ThemeColor = 1
$tr.Range("D3").Font.TintAndShade = 0.499984740745262
$tr.Range("D4").Font.ThemeColor = 1
$tr.Range("D4").Font.TintAndShade = 0.499984740745262
$tr.Range("D7").Font.ThemeColor = 1
$tr.Range("D7").Font.TintAndShade = 0.499984740745262
$tr.Range("D8").Font.ThemeColor = 1
$tr.Range("D8").Font.TintAndShade = 0.499984740745262
$tr.Range("D9").Font.ThemeColor = 1
$tr.Range("D9").Font.TintAndShade = 0.499984740745262

$tr.Range("E3").Font.ThemeColor = 1
$tr.Range("E3").Font.TintAndShade = 0.499984740745262
$tr.Range("E4").Font.ThemeColor = 1
$tr.Range("E4").Font.TintAndShade = 0.499984740745262
$tr.Range("E7").Font.ThemeColor = 1
$tr.Range("E7").Font.TintAndShade = 0.499984740745262
$tr.Range("E8").Font.ThemeColor = 1
$tr.Range("E8").Font.TintAndShade = 0.499984740745262
$tr.Range("E9").Font.ThemeColor = 1
$tr.Range("E9").Font.TintAndShade = 0.499984740745262

$tr.Range("F3").Font.Bold = $true
$tr.Range("F4").Font.Bold = $true
$tr.Range("F7").Font.Bold = $true
$tr.Range("F8").Font.Bold = $true
$tr.Range("F9").Font.Bold = $true

$tr.Range("B15").Font.ThemeColor = 1
$tr.Range("B15").Font.TintAndShade = 0.499984740745262

# Electricity row: remove stale fixed figures, explain why instead
$tr.Range("G5").ClearContents()
$tr.Range("D5").Value = "Fixed no matter production level"
$tr.Range("E5").Value = "Modeled as separate electricity consuming component"
$tr.Range("F5").ClearContents()
$tr.Range("H5").ClearContents()

$tr.Range("G19").Select()

$wb.Worksheets.Item("Capacity_Market").Range("E33").Select()

$htse.Activate()
